$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 42651.599548611113
$ws.Range("B2").Value = 14
$ws.Range("D2").Value = 44
$ws.Range("E2").Value = 8268
$ws.Range("F2").Value = 949
$ws.Range("I2").Value = 93
$ws.Range("J2").Value = 6
$ws.Range("K2").Value = 42059
$ws.Range("L2").Value = 112
$ws.Range("M2").Value = 72
$ws.Range("N2").Value = 73
$ws.Range("O2").Value = 5
$ws.Range("Q2").Value = 46.242130528922125
$ws.Range("W2").Value = 1
